$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived LR-pair data for Amelx-Cd63: now includes "ECs" as an
# additional sending cluster, giving a full 3x3 Sending x Target cluster
# matrix (ECs, FAPs, MuSCs) in place of the previous 2x3 matrix.
$data = @(
  @("ECs", "Amelx", "Cd63", "ECs", 1, 0.3333333333333333, 0.052706, 0.158118, 0.0131977395622021, 0.0131977395622021, 3, 1, 2.307850666666667, 6.923552000000001, 0.0151489512424924, 0.0151489512424924, 0.1216375772373334, 1.094738195136, 0.0001999319131389127, 0.0001999319131389127),
  @("ECs", "Amelx", "Cd63", "FAPs", 1, 0.3333333333333333, 0.052706, 0.158118, 0.0131977395622021, 0.0131977395622021, 3, 1, 110.7417803333333, 332.225341, 0.726919577166375, 0.726919577166375, 5.836756274248668, 52.530806468238, 0.009593695262107894, 0.009593695262107894),
  @("ECs", "Amelx", "Cd63", "MuSCs", 1, 0.3333333333333333, 0.052706, 0.158118, 0.0131977395622021, 0.0131977395622021, 3, 1, 39.29429233333333, 117.882877, 0.2579314715911324, 0.2579314715911324, 2.071044971720667, 18.639404745486, 0.003404112386955297, 0.003404112386955297),
  @("FAPs", "Amelx", "Cd63", "ECs", 3, 1, 1.388571333333333, 4.165713999999999, 0.347702402399595, 0.347702402399595, 3, 1, 2.307850666666667, 6.923552000000001, 0.0151489512424924, 0.0151489512424924, 3.204615277347556, 28.841537496128, 0.005267326740848937, 0.005267326740848937),
  @("FAPs", "Amelx", "Cd63", "FAPs", 3, 1, 1.388571333333333, 4.165713999999999, 0.347702402399595, 0.347702402399595, 3, 1, 110.7417803333333, 332.225341, 0.726919577166375, 0.726919577166375, 153.7728615731638, 1383.955754158474, 0.2527516833320464, 0.2527516833320464),
  @("FAPs", "Amelx", "Cd63", "MuSCs", 3, 1, 1.388571333333333, 4.165713999999999, 0.347702402399595, 0.347702402399595, 3, 1, 39.29429233333333, 117.882877, 0.2579314715911324, 0.2579314715911324, 54.56292789768644, 491.0663510791779, 0.08968339232669965, 0.08968339232669965),
  @("MuSCs", "Amelx", "Cd63", "ECs", 3, 1, 2.552285333333333, 7.656856, 0.6390998580382028, 0.6390998580382028, 3, 1, 2.307850666666667, 6.923552000000001, 0.0151489512424924, 0.0151489512424924, 5.890293408056889, 53.01264067251201, 0.009681692588504548, 0.009681692588504548),
  @("MuSCs", "Amelx", "Cd63", "FAPs", 3, 1, 2.552285333333333, 7.656856, 0.6390998580382028, 0.6390998580382028, 3, 1, 110.7417803333333, 332.225341, 0.726919577166375, 0.726919577166375, 282.6446217319885, 2543.801595587896, 0.4645741985722207, 0.4645741985722207),
  @("MuSCs", "Amelx", "Cd63", "MuSCs", 3, 1, 2.552285333333333, 7.656856, 0.6390998580382028, 0.6390998580382028, 3, 1, 39.29429233333333, 117.882877, 0.2579314715911324, 0.2579314715911324, 100.2902460060791, 902.612214054712, 0.1648439668774775, 0.1648439668774775)
)

$rowCount = $data.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $rowValues = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $rowValues[$j]
    }
}
